# Apply updated cryptocurrency price/volume data per the Nov 6 2023 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "35.560.76"
$ws.Cells.Item(2, 5).Value = "  +0.22%  "
$ws.Cells.Item(3, 4).Value = "1.913.25"
$ws.Cells.Item(3, 5).Value = "  +0.47%  "
$ws.Cells.Item(4, 5).Value = "  -0.14%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.707"
$ws.Cells.Item(5, 5).Value = "  +9.01%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "246.89"
$ws.Cells.Item(6, 5).Value = "  +0.07%  "
$ws.Cells.Item(7, 5).Value = "  -0.05%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "40.91"
$ws.Cells.Item(8, 5).Value = "  -2.43%  "
$ws.Cells.Item(9, 5).Value = "  +3.82%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "52.67"
$ws.Cells.Item(10, 5).Value = "  +7.73%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0732"
$ws.Cells.Item(11, 5).Value = "  +3.50%  "
$ws.Cells.Item(12, 5).Value = "  -1.02%  "
$ws.Cells.Item(13, 4).Value = "2.189.28"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "12.62"
$ws.Cells.Item(14, 5).Value = "  +1.79%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.717"
$ws.Cells.Item(15, 5).Value = "  +2.67%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "4.92"
$ws.Cells.Item(16, 5).Value = "  +1.89%  "
$ws.Cells.Item(17, 4).Value = "1.897.86"
$ws.Cells.Item(17, 5).Value = "  -0.45%  "
$ws.Cells.Item(18, 4).Value = "35.566.33"
$ws.Cells.Item(18, 5).Value = "  +0.30%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "73.31"
$ws.Cells.Item(19, 5).Value = "  +1.72%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0828"
$ws.Cells.Item(20, 5).Value = "  -0.31%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.17"
$ws.Cells.Item(21, 5).Value = "  +4.17%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "242.88"
$ws.Cells.Item(22, 5).Value = "  -0.13%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.04"
$ws.Cells.Item(23, 5).Value = "  +4.29%  "
$ws.Cells.Item(24, 5).Value = "  -0.15%  "
$ws.Cells.Item(25, 5).Value = "  +0.89%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.31"
$ws.Cells.Item(26, 5).Value = "  +2.57%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "169.45"
$ws.Cells.Item(27, 5).Value = "  -1.25%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "8.67"
$ws.Cells.Item(28, 5).Value = "  +1.61%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "18.80"
$ws.Cells.Item(29, 5).Value = "  +4.62%  "
$ws.Cells.Item(30, 5).Value = "  +4.15%  "
$ws.Cells.Item(31, 4).Value = "4.119.90"
$ws.Cells.Item(31, 5).Value = "  +19.20%  "
$ws.Cells.Item(32, 5).Value = "  +2.34%  "
$ws.Cells.Item(33, 5).Value = "  +0.91%  "
$ws.Cells.Item(34, 5).Value = "  +0.52%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.88"
$ws.Cells.Item(35, 5).Value = "  +6.17%  "
$ws.Cells.Item(36, 5).Value = "  -0.09%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.916"
$ws.Cells.Item(37, 5).Value = "  -4.71%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.49"
$ws.Cells.Item(38, 5).Value = "  +11.63%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.04"
$ws.Cells.Item(39, 5).Value = "  +0.69%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "17.23"
$ws.Cells.Item(40, 5).Value = "  +9.96%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "97.88"
$ws.Cells.Item(41, 5).Value = "  +6.61%  "
$ws.Cells.Item(42, 5).Value = "  +1.44%  "
$ws.Cells.Item(43, 5).Value = "  +2.80%  "
$ws.Cells.Item(44, 5).Value = "  +1.41%  "
$ws.Cells.Item(45, 4).Value = "1.357.74"
$ws.Cells.Item(45, 5).Value = "  +0.76%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.45"
$ws.Cells.Item(46, 5).Value = "  +2.30%  "
$ws.Cells.Item(47, 2).Value = "HuobiToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.42"
$ws.Cells.Item(47, 5).Value = "  +0.31%  "
$ws.Cells.Item(48, 2).Value = "MultiversX"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "46.12"
$ws.Cells.Item(48, 5).Value = "  -6.17%  "
$ws.Cells.Item(49, 5).Value = "  +1.29%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "12.21"
$ws.Cells.Item(50, 5).Value = "  -3.30%  "
$ws.Cells.Item(51, 5).Value = "  -0.80%  "
